# Updated cryptos list with latest Price / Volume(1h) figures.
# D/E columns hold text-formatted figures (e.g. "65.869.72", "  +7.13%  ");
# for D-column values that look numeric, force text format first so Excel
# doesn't silently coerce them to a Double (which would also mangle strings
# like "463.80" -> 463.8 / "8.50" -> 8.5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.869.72"
$ws.Range("E2").Value = "  +7.13%  "

$ws.Range("D3").Value = "3.015.76"
$ws.Range("E3").Value = "  +4.39%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.54%  "

$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").Value = "3.012.25"
$ws.Range("E8").Value = "  +4.37%  "

$ws.Range("E9").Value = "  +3.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.77%  "

$ws.Range("E11").Value = "  +7.94%  "

$ws.Range("E12").Value = "  +5.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.97%  "

$ws.Range("E15").Value = "  +0.63%  "

$ws.Range("D16").Value = "65.923.29"
$ws.Range("E16").Value = "  +7.20%  "

$ws.Range("D17").Value = "3.516.07"
$ws.Range("E17").Value = "  +4.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.90%  "

$ws.Range("D19").Value = "3.014.95"
$ws.Range("E19").Value = "  +4.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.65%  "

$ws.Range("E22").Value = "  +5.07%  "

$ws.Range("E23").Value = "  +8.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.71%  "

$ws.Range("E26").Value = "  +12.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.33%  "

$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +14.12%  "

$ws.Range("E30").Value = "  +18.52%  "

$ws.Range("E31").Value = "  +0.77%  "

$ws.Range("E32").Value = "  +5.28%  "

$ws.Range("E33").Value = "  +5.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.40%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.22%  "

$ws.Range("E37").Value = "  +8.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "44.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.51%  "

$ws.Range("E42").Value = "  +8.94%  "

$ws.Range("E43").Value = "  +12.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "398.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +16.64%  "

$ws.Range("D46").Value = "2.805.93"
$ws.Range("E46").Value = "  +4.53%  "

$ws.Range("E47").Value = "  +6.30%  "

$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.03%  "

$ws.Range("E51").Value = "  +4.70%  "

